# Negate all numeric values in column E (Block frame data) for rows 2-136.
# Cells that are empty / non-numeric are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 136; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -is [double]) {
        $cell.Value = (0 - $val)
    }
}
